$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "U2"  = 45700
    "U4"  = 40000
    "U5"  = 40600
    "U6"  = 43700
    "U7"  = 43400
    "U8"  = 42500
    "U9"  = 40300
    "U10" = 39600
    "U11" = 38700
    "U12" = 38900
    "U13" = 35300
    "U14" = 35400
    "U15" = 36200
    "U16" = 38200
    "U17" = 39500
    "U18" = 40000
    "U19" = 41000
    "U20" = 41800
    "U21" = 43200
    "U22" = 44000
    "U23" = 42900
    "U24" = 43400
    "U25" = 42600
    "U26" = 43900
    "U27" = 44700
    "U29" = 48400
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
